$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 412.2
$ws.Range("I33").Value = 291.44446
$ws.Range("K33").Value = 291.44446
$ws.Range("M33").Value = -62.44445999999999
$ws.Range("H76").Value = 83341740
$ws.Range("I76").Value = 166671500
$ws.Range("J76").Value = 11995.667
$ws.Range("K76").Value = 166671500
$ws.Range("L76").Value = 11995.667
$ws.Range("M76").Value = -166671185
$ws.Range("N76").Value = -12625.667
$ws.Range("H79").Value = 83341740
$ws.Range("I79").Value = 166671500
$ws.Range("J79").Value = 11995.667
$ws.Range("K79").Value = 166671500
$ws.Range("L79").Value = 11995.667
$ws.Range("M79").Value = -166670408
$ws.Range("N79").Value = -14179.667
$ws.Range("H98").Value = 1361.5518
$ws.Range("I98").Value = 693.17645
$ws.Range("J98").Value = 2308.4167
$ws.Range("K98").Value = 693.17645
$ws.Range("L98").Value = 2308.4167
$ws.Range("M98").Value = 804.82355
$ws.Range("N98").Value = -5304.4167
$ws.Range("H99").Value = 511.58334
$ws.Range("I99").Value = 252
$ws.Range("J99").Value = 697
$ws.Range("K99").Value = 756
$ws.Range("L99").Value = 2091
$ws.Range("M99").Value = 742
$ws.Range("N99").Value = -5087
$ws.Range("H122").Value = 1361.5518
$ws.Range("I122").Value = 693.17645
$ws.Range("J122").Value = 2308.4167
$ws.Range("K122").Value = 2079.52935
$ws.Range("L122").Value = 6925.250100000001
$ws.Range("M122").Value = 370.4706499999998
$ws.Range("N122").Value = -11825.2501
$ws.Range("H132").Value = 2802.65
$ws.Range("I132").Value = 2586.2778
$ws.Range("K132").Value = 7758.8334
$ws.Range("M132").Value = -5228.8334
$ws.Range("H137").Value = 310761.44
$ws.Range("I137").Value = 2191.25
$ws.Range("J137").Value = 765496.5
$ws.Range("K137").Value = 6573.75
$ws.Range("L137").Value = 2296489.5
$ws.Range("M137").Value = -4023.75
$ws.Range("N137").Value = -2301589.5
$ws.Range("H138").Value = 1620.9302
$ws.Range("I138").Value = 1330.6364
$ws.Range("K138").Value = 3991.9092
$ws.Range("M138").Value = 1148.0908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6674.4517
$ws.Range("I32").Value = 3437.0425
$ws.Range("K32").Value = 3437.0425
$ws.Range("M32").Value = -3150.0425
$ws.Range("H44").Value = 89417.39999999999
$ws.Range("J44").Value = 94271.75
$ws.Range("L44").Value = 94271.75
$ws.Range("N44").Value = -95247.75
$ws.Range("H45").Value = 10422512
$ws.Range("I45").Value = 3177.2856
$ws.Range("J45").Value = 25009580
$ws.Range("K45").Value = 3177.2856
$ws.Range("L45").Value = 25009580
$ws.Range("M45").Value = -2800.2856
$ws.Range("N45").Value = -25010334
$ws.Range("H55").Value = 50000
$ws.Range("J55").Value = 50000
$ws.Range("L55").Value = 50000
$ws.Range("N55").Value = -50630
$ws.Range("H74").Value = 50637.855
$ws.Range("I74").Value = 126862.125
$ws.Range("J74").Value = 3730.6155
$ws.Range("K74").Value = 126862.125
$ws.Range("L74").Value = 3730.6155
$ws.Range("M74").Value = -125988.125
$ws.Range("N74").Value = -5478.6155
$ws.Range("H77").Value = 50637.855
$ws.Range("I77").Value = 126862.125
$ws.Range("J77").Value = 3730.6155
$ws.Range("K77").Value = 634310.625
$ws.Range("L77").Value = 18653.0775
$ws.Range("M77").Value = -629942.625
$ws.Range("N77").Value = -27389.0775
$ws.Range("H122").Value = 3882
$ws.Range("I122").Value = 3927
$ws.Range("K122").Value = 11781
$ws.Range("M122").Value = -9331
$ws.Range("H132").Value = 2134.652
$ws.Range("I132").Value = 1474.9375
$ws.Range("J132").Value = 3642.5715
$ws.Range("K132").Value = 4424.8125
$ws.Range("L132").Value = 10927.7145
$ws.Range("M132").Value = -1894.8125
$ws.Range("N132").Value = -15987.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 336817.72
$ws.Range("I20").Value = 470064.8
$ws.Range("J20").Value = 3700
$ws.Range("K20").Value = 470064.8
$ws.Range("L20").Value = 3700
$ws.Range("M20").Value = -469817.8
$ws.Range("N20").Value = -4194
$ws.Range("H107").Value = 2107.5334
$ws.Range("I107").Value = 1920.5
$ws.Range("J107").Value = 2481.6
$ws.Range("K107").Value = 1920.5
$ws.Range("L107").Value = 2481.6
$ws.Range("M107").Value = -0.5
$ws.Range("N107").Value = -6321.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7166666.5
$ws.Range("I6").Value = 10250000
$ws.Range("K6").Value = 10250000
$ws.Range("M6").Value = -10249887
$ws.Range("H31").Value = 2111.9375
$ws.Range("I31").Value = 1726.0667
$ws.Range("K31").Value = 1726.0667
$ws.Range("M31").Value = -1431.0667
$ws.Range("H34").Value = 2111.9375
$ws.Range("I34").Value = 1726.0667
$ws.Range("K34").Value = 1726.0667
$ws.Range("M34").Value = -1524.0667
$ws.Range("H62").Value = 2818.6
$ws.Range("J62").Value = 2364.6667
$ws.Range("L62").Value = 2364.6667
$ws.Range("N62").Value = -3612.6667
$ws.Range("H65").Value = 2818.6
$ws.Range("J65").Value = 2364.6667
$ws.Range("L65").Value = 11823.3335
$ws.Range("N65").Value = -18063.3335
$ws.Range("H141").Value = 101444.5
$ws.Range("J141").Value = 114890
$ws.Range("L141").Value = 114890
$ws.Range("N141").Value = -125250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 859.8889
$ws.Range("I5").Value = 650.36365
$ws.Range("K5").Value = 1951.09095
$ws.Range("M5").Value = -1839.09095
$ws.Range("H107").Value = 397.5
$ws.Range("J107").Value = 397.5
$ws.Range("L107").Value = 1192.5
$ws.Range("N107").Value = -5032.5
$ws.Range("H113").Value = 46376.773
$ws.Range("I113").Value = 831.125
$ws.Range("K113").Value = 2493.375
$ws.Range("M113").Value = -323.375
$ws.Range("H135").Value = 859.8889
$ws.Range("I135").Value = 650.36365
$ws.Range("K135").Value = 5853.27285
$ws.Range("M135").Value = -3318.27285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7798.3
$ws.Range("J70").Value = 7975
$ws.Range("L70").Value = 7975
$ws.Range("N70").Value = -8515
$ws.Range("H73").Value = 7798.3
$ws.Range("J73").Value = 7975
$ws.Range("L73").Value = 7975
$ws.Range("N73").Value = -9847
$ws.Range("H97").Value = 1419.1666
$ws.Range("I97").Value = 693.55
$ws.Range("J97").Value = 5047.25
$ws.Range("K97").Value = 693.55
$ws.Range("L97").Value = 5047.25
$ws.Range("M97").Value = -197.55
$ws.Range("N97").Value = -6039.25
$ws.Range("H122").Value = 114331.6
$ws.Range("I122").Value = 125256.5
$ws.Range("K122").Value = 375769.5
$ws.Range("M122").Value = -373319.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3625.875
$ws.Range("I132").Value = 3167.8333
$ws.Range("K132").Value = 9503.499899999999
$ws.Range("M132").Value = -6973.499899999999
$ws.Range("H136").Value = 1921.1305
$ws.Range("I136").Value = 1723.1428
$ws.Range("K136").Value = 5169.428400000001
$ws.Range("M136").Value = -2619.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13631
$ws.Range("I132").Value = 17548.053
$ws.Range("K132").Value = 52644.159
$ws.Range("M132").Value = -50114.159
